$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.082653666666667
$ws.Range("H2").Value = 6.247961
$ws.Range("I2").Value = 0.0472190032704503
$ws.Range("J2").Value = 0.0472190032704503
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 4.171693666666666
$ws.Range("N2").Value = 12.515081
$ws.Range("O2").Value = 0.1077921033402881
$ws.Range("P2").Value = 0.1077921033402881
$ws.Range("Q2").Value = 8.688193111093444
$ws.Range("R2").Value = 78.19373799984099
$ws.Range("S2").Value = 0.005089835680153781
$ws.Range("T2").Value = 0.005089835680153781
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.082653666666667
$ws.Range("H3").Value = 6.247961
$ws.Range("I3").Value = 0.0472190032704503
$ws.Range("J3").Value = 0.0472190032704503
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 20.39394566666667
$ws.Range("N3").Value = 61.181837
$ws.Range("O3").Value = 0.5269577477327285
$ws.Range("P3").Value = 0.5269577477327286
$ws.Range("Q3").Value = 42.47352572048412
$ws.Range("R3").Value = 382.261731484357
$ws.Range("S3").Value = 0.02488241961358083
$ws.Range("T3").Value = 0.02488241961358083
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.082653666666667
$ws.Range("H4").Value = 6.247961
$ws.Range("I4").Value = 0.0472190032704503
$ws.Range("J4").Value = 0.0472190032704503
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 14.13565266666667
$ws.Range("N4").Value = 42.406958
$ws.Range("O4").Value = 0.3652501489269833
$ws.Range("P4").Value = 0.3652501489269833
$ws.Range("Q4").Value = 29.43966885695978
$ws.Range("R4").Value = 264.957019712638
$ws.Range("S4").Value = 0.01724674797671568
$ws.Range("T4").Value = 0.01724674797671568
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 37.28222
$ws.Range("H5").Value = 111.84666
$ws.Range("I5").Value = 0.8452818134314446
$ws.Range("J5").Value = 0.8452818134314446
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.171693666666666
$ws.Range("N5").Value = 12.515081
$ws.Range("O5").Value = 0.1077921033402881
$ws.Range("P5").Value = 0.1077921033402881
$ws.Range("Q5").Value = 155.5300010532733
$ws.Range("R5").Value = 1399.77000947946
$ws.Range("S5").Value = 0.09111470458506842
$ws.Range("T5").Value = 0.09111470458506843
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 37.28222
$ws.Range("H6").Value = 111.84666
$ws.Range("I6").Value = 0.8452818134314446
$ws.Range("J6").Value = 0.8452818134314446
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 20.39394566666667
$ws.Range("N6").Value = 61.181837
$ws.Range("O6").Value = 0.5269577477327285
$ws.Range("P6").Value = 0.5269577477327286
$ws.Range("Q6").Value = 760.3315690127132
$ws.Range("R6").Value = 6842.984121114419
$ws.Range("S6").Value = 0.4454278006052705
$ws.Range("T6").Value = 0.4454278006052706
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 37.28222
$ws.Range("H7").Value = 111.84666
$ws.Range("I7").Value = 0.8452818134314446
$ws.Range("J7").Value = 0.8452818134314446
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 14.13565266666667
$ws.Range("N7").Value = 42.406958
$ws.Range("O7").Value = 0.3652501489269833
$ws.Range("P7").Value = 0.3652501489269833
$ws.Range("Q7").Value = 527.0085125622533
$ws.Range("R7").Value = 4743.07661306028
$ws.Range("S7").Value = 0.3087393082411056
$ws.Range("T7").Value = 0.3087393082411057
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 4.741387
$ws.Range("H8").Value = 14.224161
$ws.Range("I8").Value = 0.107499183298105
$ws.Range("J8").Value = 0.107499183298105
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 4.171693666666666
$ws.Range("N8").Value = 12.515081
$ws.Range("O8").Value = 0.1077921033402881
$ws.Range("P8").Value = 0.1077921033402881
$ws.Range("Q8").Value = 19.77961411911567
$ws.Range("R8").Value = 178.016527072041
$ws.Range("S8").Value = 0.01158756307506591
$ws.Range("T8").Value = 0.01158756307506591
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 4.741387
$ws.Range("H9").Value = 14.224161
$ws.Range("I9").Value = 0.107499183298105
$ws.Range("J9").Value = 0.107499183298105
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 20.39394566666667
$ws.Range("N9").Value = 61.181837
$ws.Range("O9").Value = 0.5269577477327285
$ws.Range("P9").Value = 0.5269577477327286
$ws.Range("Q9").Value = 96.69558886263968
$ws.Range("R9").Value = 870.2602997637571
$ws.Range("S9").Value = 0.05664752751387718
$ws.Range("T9").Value = 0.05664752751387718
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 4.741387
$ws.Range("H10").Value = 14.224161
$ws.Range("I10").Value = 0.107499183298105
$ws.Range("J10").Value = 0.107499183298105
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 14.13565266666667
$ws.Range("N10").Value = 42.406958
$ws.Range("O10").Value = 0.3652501489269833
$ws.Range("P10").Value = 0.3652501489269833
$ws.Range("Q10").Value = 67.02259979024868
$ws.Range("R10").Value = 603.2033981122381
$ws.Range("S10").Value = 0.03926409270916194
$ws.Range("T10").Value = 0.03926409270916194
